{"js": "// Translate the bold field-labels of the use-case template from French to\n// English (template localisation). Every label lives in its own run\n// (bold/bCs) inside its own paragraph, so each one is located with\n// Paragraph.search (exact, case-sensitive) and replaced in place with\n// Range.insertText(..., \"Replace\") - this keeps the surrounding run\n// formatting (bold) untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Simple whole-paragraph-label replacements: [paragraph text prefix used to\n// find the right paragraph unambiguously, exact text to search for inside\n// it, replacement text].\nconst simpleReplacements = [\n  [\"Nom\", \"Nom\", \"Name\"],\n  [\"Identifiant\", \"Identifiant\", \"Identifier\"],\n  [\"Acteurs\", \"Acteurs\", \"Actors\"],\n  [\"Pr\u00e9conditions\", \"Pr\u00e9conditions\", \"Preconditions\"],\n  [\"Hypoth\u00e8ses\", \"Hypoth\u00e8ses\", \"Hypothesis\"],\n  [\"Fr\u00e9quence\", \"Fr\u00e9quence\", \"Frequency\"],\n  [\"Sc\u00e9nario nominal\", \"Sc\u00e9nario nominal\", \"Nominal script\"],\n  [\"Sc\u00e9nario(s) alternatif(s)\", \"Sc\u00e9nario(s) alternatif(s)\", \"Alternative script\"],\n  [\"Cas A : une description br\u00e8ve du cas\", \"Cas A : une description br\u00e8ve du cas\", \"Case A : \"],\n  [\"Cas B : une description br\u00e8ve du cas\", \"Cas B : une description br\u00e8ve du cas\", \"Case B : \"],\n  [\"Cas d\\u2019extension\", \"Cas d\\u2019extension\", \"Extension cases\"],\n];\n\nfor (const [paraPrefix, search, replacement] of simpleReplacements) {\n  const paragraph = paragraphs.items.find((p) => p.text.indexOf(paraPrefix) !== -1);\n  if (!paragraph) continue;\n  const found = paragraph.search(search, { matchCase: true });\n  found.load(\"text\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(replacement, \"Replace\");\n  }\n  await context.sync();\n}\n\n// \"Cas d'inclusion : \" is special: the label becomes \"Inclusion cases\" AND\n// the plain run right after it loses its leading space (\" : \" -> \": \").\nconst inclusionParagraph = paragraphs.items.find((p) => p.text.indexOf(\"Cas d\\u2019inclusion\") !== -1);\nif (inclusionParagraph) {\n  const label = inclusionParagraph.search(\"Cas d\\u2019inclusion\", { matchCase: true });\n  const colon = inclusionParagraph.search(\" : \", { matchCase: true });\n  label.load(\"text\");\n  colon.load(\"text\");\n  await context.sync();\n\n  if (label.items.length > 0) {\n    label.items[0].insertText(\"Inclusion cases\", \"Replace\");\n  }\n  if (colon.items.length > 0) {\n    colon.items[0].insertText(\": \", \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Translate the bold field-labels of the use-case template from French to\n# English. Each label is a single bold run living in its own paragraph, so\n# every replacement is scoped with Find.Execute() to the owning\n# Paragraph.Range - this both disambiguates short labels (e.g. \"Nom\") and\n# keeps the surrounding run's bold formatting untouched.\n\n$d = $word.ActiveDocument\n$paragraphs = $d.Paragraphs\n\nfunction Replace-InParagraph($paragraphIndex, $searchText, $replaceText) {\n    $range = $paragraphs.Item($paragraphIndex).Range\n    $range.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\nReplace-InParagraph 1  \"Nom\"                                   \"Name\"\nReplace-InParagraph 2  \"Identifiant\"                            \"Identifier\"\nReplace-InParagraph 3  \"Acteurs\"                                \"Actors\"\nReplace-InParagraph 5  \"Pr\u00e9conditions\"                          \"Preconditions\"\nReplace-InParagraph 8  \"Hypoth\u00e8ses\"                             \"Hypothesis\"\nReplace-InParagraph 11 \"Fr\u00e9quence\"                              \"Frequency\"\nReplace-InParagraph 12 \"Sc\u00e9nario nominal\"                       \"Nominal script\"\nReplace-InParagraph 16 \"Sc\u00e9nario(s) alternatif(s)\"              \"Alternative script\"\nReplace-InParagraph 17 \"Cas A : une description br\u00e8ve du cas\"   \"Case A : \"\nReplace-InParagraph 20 \"Cas B : une description br\u00e8ve du cas\"   \"Case B : \"\nReplace-InParagraph 25 \"Cas d\u2019extension\"                        \"Extension cases\"\n\n# \"Cas d'inclusion : \" needs two changes in its paragraph: the label becomes\n# \"Inclusion cases\" AND the following plain run loses its leading space\n# (\" : \" -> \": \").\nReplace-InParagraph 24 \"Cas d\u2019inclusion\" \"Inclusion cases\"\nReplace-InParagraph 24 \" : \"            \": \"\n"}
